$d = $word.ActiveDocument

# 1. Replace the ID placeholder text (and remove the trailing space run) in the
#    first paragraph with the new placeholder text.
$d.Content.Find.Execute("**ID__AFFARS_5317_topic_2__ID** ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_SUBPART_5317_1__ID**", 2)

# 2. Update the first paragraph's formatting: add a paragraph border (5 twips
#    space on every side, no explicit line) and change the left indent from
#    120 twips (6pt) to 225 twips (11.25pt).
$p1 = $d.Paragraphs.Item(1)
$pf = $p1.Range.ParagraphFormat
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5
$pf.LeftIndent = 11.25
